# "Ajuste el tema costos" -- update Sprint 2's hours in the
# "Cobertura de la Prueba" sheet: the planned-hours formula now adds 76
# (instead of 74) for this sprint, and the completed / actual hours for
# the sprint are corrected. Rows 5:7 hold plain "+74" running-total
# formulas, so they recalculate automatically from the new B4 value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cobertura de la Prueba")

# Hs Planificadas Totales (Sprint 2): formula bumped from B3+74 to B3+76.
$ws.Range("B4").Formula = "=B3+76"

# Hs Planificadas Completadas (Sprint 2): 130 -> 144.
$ws.Range("C4").Value = 144

# Hs Reales Utilizadas (Sprint 2): 145 -> 152.
$ws.Range("D4").Value = 152

# Recalculate so the dependent running totals (B5:B7) pick up the change.
$excel.CalculateFullRebuild()

# Leave the selection on the last-edited cell, like Excel would after
# typing the new value and pressing Enter.
$ws.Range("D4").Select()
